$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.5341408786989064
$ws.Range("C2").Value = 0.03856766815198398
$ws.Range("D2").Value = 0.617020563227868
$ws.Range("E2").Value = 0.2351380644541905
$ws.Range("G2").Value = 0.002515417273624479
$ws.Range("J2").Value = 0.1101405180386799
$ws.Range("K2").Value = 0.4952642386334389
$ws.Range("M2").Value = 0.3580431246268532
$ws.Range("N2").Value = 2.38255128027653
$ws.Range("O2").Value = 5.534750650312361
$ws.Range("B3").Value = 0.5015122491333841
$ws.Range("C3").Value = 0.03446157335696398
$ws.Range("D3").Value = 0.6116952897218084
$ws.Range("E3").Value = 0.2337616974073136
$ws.Range("G3").Value = 0.002518313550613195
$ws.Range("J3").Value = 0.1099336660209254
$ws.Range("K3").Value = 0.4611112161629762
$ws.Range("M3").Value = 0.347266301744078
$ws.Range("N3").Value = 2.403493213796065
$ws.Range("O3").Value = 5.548900556750908
$ws.Range("B4").Value = 0.4817097504417234
$ws.Range("C4").Value = 0.03194520565540415
$ws.Range("D4").Value = 0.6087097194675692
$ws.Range("E4").Value = 0.2330216455758141
$ws.Range("G4").Value = 0.002520187734710964
$ws.Range("J4").Value = 0.1098528224466619
$ws.Range("K4").Value = 0.4403359474610511
$ws.Range("M4").Value = 0.3408242104516006
$ws.Range("N4").Value = 2.417037316283132
$ws.Range("O4").Value = 5.560034520731392
$ws.Range("B5").Value = 0.4736986863694312
$ws.Range("C5").Value = 0.03092100569445222
$ws.Range("D5").Value = 0.6075646538674704
$ws.Range("E5").Value = 0.2327465199450494
$ws.Range("G5").Value = 0.002520975657303036
$ws.Range("J5").Value = 0.1098314981509994
$ws.Range("K5").Value = 0.431919215132524
$ws.Range("M5").Value = 0.3382431284478997
$ws.Range("N5").Value = 2.422729072984581
$ws.Range("O5").Value = 5.565186979899124
$ws.Range("B6").Value = 0.4723720065261716
$ws.Range("C6").Value = 0.03075101403499048
$ws.Range("D6").Value = 0.6073788433634348
$ws.Range("E6").Value = 0.232702434259501
$ws.Range("G6").Value = 0.002521107953705346
$ws.Range("J6").Value = 0.1098286594410283
$ws.Range("K6").Value = 0.4305246151293147
$ws.Range("M6").Value = 0.3378172106862394
$ws.Range("N6").Value = 2.423684599160424
$ws.Range("O6").Value = 5.566079710481461
$ws.Range("B7").Value = 0.4816014726041828
$ws.Range("C7").Value = 0.03193138786761551
$ws.Range("D7").Value = 0.6086939867307564
$ws.Range("E7").Value = 0.2330178279840212
$ws.Range("G7").Value = 0.00252019826297553
$ws.Range("J7").Value = 0.1098524877949032
$ws.Range("K7").Value = 0.4402222361269992
$ws.Range("M7").Value = 0.3407892222444957
$ws.Range("N7").Value = 2.417113379155786
$ws.Range("O7").Value = 5.560101517041261
$ws.Range("B8").Value = 0.5228426694909558
$ws.Range("C8").Value = 0.03715091228846745
$ws.Range("D8").Value = 0.6151254995233018
$ws.Range("E8").Value = 0.2346417129282869
$ws.Range("G8").Value = 0.002516396058402107
$ws.Range("J8").Value = 0.1100596207356723
$ws.Range("K8").Value = 0.4834480497289348
$ws.Range("M8").Value = 0.3542910491137974
$ws.Range("N8").Value = 2.389629754294383
$ws.Range("O8").Value = 5.539122036687075
$ws.Range("B9").Value = 0.6055423739819048
$ws.Range("C9").Value = 0.04742335792768415
$ws.Range("D9").Value = 0.6299884206613626
$ws.Range("E9").Value = 0.2386583919828666
$ws.Range("G9").Value = 0.002509697221063629
$ws.Range("J9").Value = 0.110831727479038
$ws.Range("K9").Value = 0.5697480821175418
$ws.Range("M9").Value = 0.3821517382726256
$ws.Range("N9").Value = 2.341173851912316
$ws.Range("O9").Value = 5.517381795879857
$ws.Range("B10").Value = 0.6674056650519447
$ws.Range("C10").Value = 0.05499243404422316
$ws.Range("D10").Value = 0.6422769177200394
$ws.Range("E10").Value = 0.2421157899573743
$ws.Range("G10").Value = 0.002505232589981769
$ws.Range("J10").Value = 0.1116217733426197
$ws.Range("K10").Value = 0.6340795092157521
$ws.Range("M10").Value = 0.4034613279940302
$ws.Range("N10").Value = 2.308886732478509
$ws.Range("O10").Value = 5.513234282034063
$ws.Range("B11").Value = 0.695787011169557
$ws.Range("C11").Value = 0.05844047531542174
$ws.Range("D11").Value = 0.6481638862585442
$ws.Range("E11").Value = 0.2437984256657515
$ws.Range("G11").Value = 0.002503299759635233
$ws.Range("J11").Value = 0.1120295127327253
$ws.Range("K11").Value = 0.6635453911624722
$ws.Range("M11").Value = 0.4133375762011866
$ws.Range("N11").Value = 2.29491743197427
$ws.Range("O11").Value = 5.513915304680268
$ws.Range("B12").Value = 0.706568434289494
$ws.Range("C12").Value = 0.0597468286703986
$ws.Range("D12").Value = 0.6504357273703363
$ws.Range("E12").Value = 0.2444513650574436
$ws.Range("G12").Value = 0.002502581886334966
$ws.Range("J12").Value = 0.1121908569908356
$ws.Range("K12").Value = 0.6747320111436181
$ws.Range("M12").Value = 0.4171035849029394
$ws.Range("N12").Value = 2.289730912365044
$ws.Range("O12").Value = 5.514542342866378
$ws.Range("B13").Value = 0.7042449551802008
$ws.Range("C13").Value = 0.05946545360271216
$ws.Range("D13").Value = 0.6499445540921442
$ws.Range("E13").Value = 0.244310042412458
$ws.Range("G13").Value = 0.00250273586942633
$ws.Range("J13").Value = 0.1121558000135536
$ws.Range("K13").Value = 0.6723215098866149
$ws.Range("M13").Value = 0.4162913495483807
$ws.Range("N13").Value = 2.290843324518065
$ws.Range("O13").Value = 5.514390882320669
$ws.Range("B14").Value = 0.6966733257840758
$ws.Range("C14").Value = 0.05854793689259452
$ws.Range("D14").Value = 0.6483499394117018
$ws.Range("E14").Value = 0.2438518276890491
$ws.Range("G14").Value = 0.002503240418460794
$ws.Range("J14").Value = 0.1120426475444418
$ws.Range("K14").Value = 0.664465151636108
$ws.Range("M14").Value = 0.413646885979432
$ws.Range("N14").Value = 2.294488662188172
$ws.Range("O14").Value = 5.513959493678271
$ws.Range("B15").Value = 0.6920399029500004
$ws.Range("C15").Value = 0.05798601553365756
$ws.Range("D15").Value = 0.6473787328258425
$ws.Range("E15").Value = 0.2435732095726735
$ws.Range("G15").Value = 0.002503551296780155
$ws.Range("J15").Value = 0.1119742421870313
$ws.Range("K15").Value = 0.6596566088963414
$ws.Range("M15").Value = 0.4120304688703982
$ws.Range("N15").Value = 2.296734998135221
$ws.Range("O15").Value = 5.513743327278178
$ws.Range("B16").Value = 0.6655556604952153
$ws.Range("C16").Value = 0.05476719024694887
$ws.Range("D16").Value = 0.641898154371745
$ws.Range("E16").Value = 0.2420080332504142
$ws.Range("G16").Value = 0.002505360874632176
$ws.Range("J16").Value = 0.1115960985932176
$ws.Range("K16").Value = 0.632157858741806
$ws.Range("M16").Value = 0.4028195503130547
$ws.Range("N16").Value = 2.309814118042464
$ws.Range("O16").Value = 5.513241437528507
$ws.Range("B17").Value = 0.6493694693526493
$ws.Range("C17").Value = 0.05279375415959464
$ws.Range("D17").Value = 0.6386119501291603
$ws.Range("E17").Value = 0.2410759585938962
$ws.Range("G17").Value = 0.002506496083701606
$ws.Range("J17").Value = 0.1113764948804601
$ws.Range("K17").Value = 0.6153395109830626
$ws.Range("M17").Value = 0.3972155681395293
$ws.Range("N17").Value = 2.31802174615828
$ws.Range("O17").Value = 5.51359114145697
$ws.Range("B18").Value = 0.6400821603568829
$ws.Range("C18").Value = 0.05165914296010499
$ws.Range("D18").Value = 0.6367497620410063
$ws.Range("E18").Value = 0.2405501953631415
$ws.Range("G18").Value = 0.002507158268356789
$ws.Range("J18").Value = 0.1112547355249021
$ws.Range("K18").Value = 0.605685003410116
$ws.Range("M18").Value = 0.394009483514921
$ws.Range("N18").Value = 2.322810159839712
$ws.Range("O18").Value = 5.514033975509875
$ws.Range("B19").Value = 0.6369415241783258
$ws.Range("C19").Value = 0.05127506280902594
$ws.Range("D19").Value = 0.63612406105932
$ws.Range("E19").Value = 0.2403739581674884
$ws.Range("G19").Value = 0.002507384062203633
$ws.Range("J19").Value = 0.1112142918113008
$ws.Range("K19").Value = 0.602419421873833
$ws.Range("M19").Value = 0.3929269129738202
$ws.Range("N19").Value = 2.32444304341918
$ws.Range("O19").Value = 5.514225425419767
$ws.Range("B20").Value = 0.6510901850912205
$ws.Range("C20").Value = 0.05300378278572282
$ws.Range("D20").Value = 0.6389588803330639
$ws.Range("E20").Value = 0.2411741094786422
$ws.Range("G20").Value = 0.002506374282512476
$ws.Range("J20").Value = 0.1113994011397281
$ws.Range("K20").Value = 0.6171278933051951
$ws.Range("M20").Value = 0.3978103448126049
$ws.Range("N20").Value = 2.317141033005871
$ws.Range("O20").Value = 5.513528901142507
$ws.Range("B21").Value = 0.6988963767276175
$ws.Range("C21").Value = 0.05881741633143633
$ws.Range("D21").Value = 0.64881716187503
$ws.Range("E21").Value = 0.2439859889632459
$ws.Range("G21").Value = 0.00250309183938995
$ws.Range("J21").Value = 0.1120756948337984
$ws.Range("K21").Value = 0.6667719833976378
$ws.Range("M21").Value = 0.414422922149825
$ws.Range("N21").Value = 2.293415132693156
$ws.Range("O21").Value = 5.514076185055046
$ws.Range("B22").Value = 0.7303386120873938
$ws.Range("C22").Value = 0.06262076119105586
$ws.Range("D22").Value = 0.6555082184259788
$ws.Range("E22").Value = 0.2459155705655292
$ws.Range("G22").Value = 0.002501028423479898
$ws.Range("J22").Value = 0.1125581519300027
$ws.Range("K22").Value = 0.6993833745052029
$ws.Range("M22").Value = 0.4254321957047011
$ws.Range("N22").Value = 2.278511328206562
$ws.Range("O22").Value = 5.516585514900328
$ws.Range("B23").Value = 0.7135393048853587
$ws.Range("C23").Value = 0.06059050990985781
$ws.Range("D23").Value = 0.6519144126328342
$ws.Range("E23").Value = 0.2448773230287884
$ws.Range("G23").Value = 0.002502122239940406
$ws.Range("J23").Value = 0.11229695653509
$ws.Range("K23").Value = 0.6819629949683588
$ws.Range("M23").Value = 0.4195424787765418
$ws.Range("N23").Value = 2.286410623373754
$ws.Range("O23").Value = 5.515049390915806
$ws.Range("B24").Value = 0.6503121923236677
$ws.Range("C24").Value = 0.05290882901603311
$ws.Range("D24").Value = 0.6388019487847032
$ws.Range("E24").Value = 0.2411297040045142
$ws.Range("G24").Value = 0.002506429319157566
$ws.Range("J24").Value = 0.1113890312214849
$ws.Range("K24").Value = 0.6163193201984996
$ws.Range("M24").Value = 0.3975413973044439
$ws.Range("N24").Value = 2.317538986185159
$ws.Range("O24").Value = 5.513556286820801
$ws.Range("B25").Value = 0.5829754020194002
$ws.Range("C25").Value = 0.04464049643574697
$ws.Range("D25").Value = 0.6257270019353314
$ws.Range("E25").Value = 0.2374827816397875
$ws.Range("G25").Value = 0.002511428845048925
$ws.Range("J25").Value = 0.1105837087441444
$ws.Range("K25").Value = 0.546238263262012
$ws.Range("M25").Value = 0.3744668825562074
$ws.Range("N25").Value = 2.353700410911578
$ws.Range("O25").Value = 5.52118637785324
